$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B/C text swaps (rows 15/16 and 40/41) ---
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'

# --- D column (Price) updates ---
# Rows whose new value is plain text (contains >1 dot, not numeric) can be set directly.
$ws.Range("D2").Value = '24.379.31'
$ws.Range("D3").Value = '1.632.81'
$ws.Range("D15").Value = '1.633.18'
$ws.Range("D24").Value = '24.372.90'
$ws.Range("D29").Value = '1.820.88'

# Rows whose new value LOOKS like a plain number (e.g. "1.011") must be forced to text,
# otherwise COM auto-converts it to a float (losing the original formatted digits/precision).
# Force text via NumberFormat "@", assign, then restore the default "Normal" style so the
# saved cell keeps the same (absent/default) style index as before the edit.
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '1.011'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '1.006'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '305.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.3619'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '47.35'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.3205'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '1.095'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.06855'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '1.009'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '5.857'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '19.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '6.484'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '0.00001035'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.06522'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '76.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '15.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '5.835'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '11.80'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '2.416'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.350'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '145.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '18.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '123.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '1.056'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '4.088'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '5.506'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '0.08353'
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '12.21'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '5.035'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.05967'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.02201'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '1.193'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '8.075'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.2015'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '1.006'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.5819'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '3.722'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '12.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.5512'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '121.23'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '1.901'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '0.06898'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '72.68'
$ws.Range("D51").Style = "Normal"

# --- E column (Volume 1h %) updates: always plain text, safe via direct .Value ---
$ws.Range("E2").Value = '  -5.90%  '
$ws.Range("E3").Value = '  -7.60%  '
$ws.Range("E4").Value = '  +0.86%  '
$ws.Range("E5").Value = '  +0.91%  '
$ws.Range("E6").Value = '  -3.72%  '
$ws.Range("E7").Value = '  -5.68%  '
$ws.Range("E8").Value = '  -6.49%  '
$ws.Range("E9").Value = '  -11.41%  '
$ws.Range("E10").Value = '  -11.28%  '
$ws.Range("E11").Value = '  -10.95%  '
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("E13").Value = '  -9.94%  '
$ws.Range("E14").Value = '  -12.60%  '
$ws.Range("E15").Value = '  -7.46%  '
$ws.Range("E16").Value = '  -8.76%  '
$ws.Range("E17").Value = '  -10.82%  '
$ws.Range("E18").Value = '  -3.93%  '
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("E20").Value = '  -12.25%  '
$ws.Range("E21").Value = '  -11.84%  '
$ws.Range("E22").Value = '  -10.69%  '
$ws.Range("E23").Value = '  -7.82%  '
$ws.Range("E24").Value = '  -5.75%  '
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("E26").Value = '  -20.62%  '
$ws.Range("E27").Value = '  -7.06%  '
$ws.Range("E28").Value = '  -10.41%  '
$ws.Range("E29").Value = '  -7.16%  '
$ws.Range("E30").Value = '  -8.15%  '
$ws.Range("E31").Value = '  -13.74%  '
$ws.Range("E32").Value = '  -4.05%  '
$ws.Range("E33").Value = '  -24.23%  '
$ws.Range("E34").Value = '  -4.63%  '
$ws.Range("E35").Value = '  -7.37%  '
$ws.Range("E36").Value = '  -14.56%  '
$ws.Range("E37").Value = '  -12.25%  '
$ws.Range("E38").Value = '  -12.13%  '
$ws.Range("E39").Value = '  -12.08%  '
$ws.Range("E40").Value = '  -7.64%  '
$ws.Range("E41").Value = '  -13.72%  '
$ws.Range("E42").Value = '  -10.74%  '
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("E44").Value = '  -11.69%  '
$ws.Range("E45").Value = '  -4.92%  '
$ws.Range("E46").Value = '  -12.88%  '
$ws.Range("E47").Value = '  -13.37%  '
$ws.Range("E48").Value = '  -8.91%  '
$ws.Range("E49").Value = '  -12.64%  '
$ws.Range("E50").Value = '  -8.16%  '
$ws.Range("E51").Value = '  -10.45%  '
